# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets to
# match the newly published scrape (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 37888
$ws.Range("G2").Value = "不可售"
$ws.Range("F5").Value = 790
$ws.Range("F11").Value = 742
$ws.Range("F12").Value = 576
$ws.Range("F13").Value = 76
$ws.Range("F15").Value = 31
$ws.Range("F16").Value = 680
$ws.Range("F20").Value = 1188
$ws.Range("G21").Value = 36
$ws.Range("F22").Value = 865
$ws.Range("F23").Value = 2580
$ws.Range("F24").Value = 1067
$ws.Range("F25").Value = 580
$ws.Range("F27").Value = 1174
$ws.Range("F29").Value = 821
$ws.Range("F30").Value = 73
$ws.Range("F31").Value = 1175

# --- 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 432
$ws.Range("F4").Value = 337

# --- 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 659

# --- 全部类型 (All types - merged view) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 659
$ws.Range("F3").Value = 37888
$ws.Range("G3").Value = "不可售"
$ws.Range("F6").Value = 790
$ws.Range("F11").Value = 432
$ws.Range("F12").Value = 337
$ws.Range("F17").Value = 742
$ws.Range("F18").Value = 576
$ws.Range("F19").Value = 76
$ws.Range("F25").Value = 31
$ws.Range("F27").Value = 680
$ws.Range("F31").Value = 1188
$ws.Range("G32").Value = 36
$ws.Range("F33").Value = 865
$ws.Range("F34").Value = 2580
$ws.Range("F35").Value = 1067
$ws.Range("F36").Value = 580
$ws.Range("F38").Value = 1174
$ws.Range("F41").Value = 821
$ws.Range("F42").Value = 73
$ws.Range("F43").Value = 1175
